$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at 214:215, pushing the existing rows 214:337 down to 216:339.
# This preserves formatting (e.g. the date style on column D) from the row that was at 214.
$ws.Range("A214:R215").Insert()

# Fill in the new row 214 (Primera) with the new weekly data point.
$ws.Range("A214").Value = 8
$ws.Range("B214").Value = "Terminal La Palmera de La Serena"
$ws.Range("C214").Value = "Coquimbo"
$ws.Range("D214").Value = 44824
$ws.Range("E214").Value = 4
$ws.Range("F214").Value = 100114014
$ws.Range("G214").Value = "Betarraga"
$ws.Range("H214").Value = "Sin especificar"
$ws.Range("I214").Value = "Primera"
$ws.Range("J214").Value = 2000
$ws.Range("K214").Value = 550
$ws.Range("L214").Value = 600
$ws.Range("M214").Value = 575
$ws.Range("N214").Value = "`$/paquete 3 unidades"
$ws.Range("O214").Value = "Provincia del Elquí"
$ws.Range("P214").Value = 192
$ws.Range("Q214").Value = 3
$ws.Range("R214").Value = "Hortaliza"

# Fill in the new row 215 (Segunda) with the new weekly data point.
$ws.Range("A215").Value = 8
$ws.Range("B215").Value = "Terminal La Palmera de La Serena"
$ws.Range("C215").Value = "Coquimbo"
$ws.Range("D215").Value = 44824
$ws.Range("E215").Value = 4
$ws.Range("F215").Value = 100114014
$ws.Range("G215").Value = "Betarraga"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Segunda"
$ws.Range("J215").Value = 1460
$ws.Range("K215").Value = 450
$ws.Range("L215").Value = 500
$ws.Range("M215").Value = 475
$ws.Range("N215").Value = "`$/paquete 3 unidades"
$ws.Range("O215").Value = "Provincia del Elquí"
$ws.Range("P215").Value = 158
$ws.Range("Q215").Value = 3
$ws.Range("R215").Value = "Hortaliza"
